$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 368 (shifts existing rows 368..380 down to 369..381,
# and carries the existing per-column formatting, e.g. the date style on column D).
$ws.Rows(368).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A368").Value = 4
$ws.Range("B368").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C368").Value = "Los Lagos"
$ws.Range("D368").Value = 45075
$ws.Range("E368").Value = 10
$ws.Range("F368").Value = "Fruta"
$ws.Range("G368").Value = 100108
$ws.Range("H368").Value = "Tropicales y subtropicales"
$ws.Range("I368").Value = 100108002
$ws.Range("J368").Value = "Mango"
$ws.Range("K368").Value = "Sin especificar"
$ws.Range("L368").Value = "Primera"
$ws.Range("M368").Value = 80
$ws.Range("N368").Value = 8500
$ws.Range("O368").Value = 9000
$ws.Range("P368").Value = 8750
$ws.Range("Q368").Value = "$/bandeja 4 kilos"
$ws.Range("R368").Value = "Perú"
$ws.Range("S368").Value = 2188
$ws.Range("T368").Value = 4
